$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.591135333333333
$ws.Range("H2").Value = 4.773406
$ws.Range("I2").Value = 0.4707829234247397
$ws.Range("J2").Value = 0.4707829234247397
$ws.Range("Q2").Value = 0.04283283279488888
$ws.Range("R2").Value = 0.3854954951539999
$ws.Range("S2").Value = 0.4707829234247397
$ws.Range("T2").Value = 0.4707829234247397

# Row 3 updates
$ws.Range("I3").Value = 0.3035973020998604
$ws.Range("J3").Value = 0.3035973020998604
$ws.Range("S3").Value = 0.3035973020998604
$ws.Range("T3").Value = 0.3035973020998604

# Row 4 updates
$ws.Range("G4").Value = 0.730693
$ws.Range("H4").Value = 2.192079
$ws.Range("I4").Value = 0.2161964349979826
$ws.Range("J4").Value = 0.2161964349979826
$ws.Range("Q4").Value = 0.01967001199566666
$ws.Range("R4").Value = 0.177030107961
$ws.Range("S4").Value = 0.2161964349979826
$ws.Range("T4").Value = 0.2161964349979826

# Row 5 updates
$ws.Range("G5").Value = 0.03184866666666667
$ws.Range("H5").Value = 0.09554600000000001
$ws.Range("I5").Value = 0.009423339477417213
$ws.Range("J5").Value = 0.009423339477417213
$ws.Range("Q5").Value = 0.0008573554904444445
$ws.Range("R5").Value = 0.007716199414
$ws.Range("S5").Value = 0.009423339477417213
$ws.Range("T5").Value = 0.009423339477417213
